$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lucene")

# New configuration labels (shared strings will be appended automatically
# when new text values are written into cells).
$rusSmote = "CV + tfidf + ngram(3) + RUS + SMOTE"
$smoteRus = "CV + tfidf + ngram(3) + SMOTE + RUS"

# Model names, in the same order used throughout the sheet.
$logReg = "Logistic Regression"
$naiveBayes = "Multinomial Naive Bayes"
$svm = "Support Vector Machines"
$decisionTree = "Decision Tree"
$randomForest = "Random Forest"

# Rows 37-41: "CV + tfidf + ngram(3) + RUS + SMOTE"
$data1 = @(
    @($logReg,      88.48, 82.91, 84.93, 92.51),
    @($naiveBayes,  86.61, 80.94, 86.66, 86.79),
    @($svm,         87.89, 80.41, 78.69, 99.89),
    @($decisionTree,69.54, 62.62, 83.8,  59.74),
    @($randomForest,79.4,  72.43, 86.34, 73.87)
)

# Rows 42-46: "CV + tfidf + ngram(3) + SMOTE + RUS"
$data2 = @(
    @($logReg,      88.52, 82.74, 84.12, 93.56),
    @($naiveBayes,  88.52, 82.99, 84.68, 92.82),
    @($svm,         87.93, 80.5,  78.84, 99.73),
    @($decisionTree,81.86, 74.07, 82.52, 81.49),
    @($randomForest,87.26, 80.78, 83.06, 92.15)
)

$row = 37
foreach ($r in $data1) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $rusSmote
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $row++
}

foreach ($r in $data2) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $smoteRus
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $row++
}

# Update the selection to match the final state in the diff.
$ws.Range("F45").Select()
